$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B11").Formula = "=SUM(B1:B10)"
$ws.Range("B11").Font.Bold = $true

$ws.Range("A1:B11").Select()

$ws.PageSetup.Orientation = 1
